# Add a new "2022-Q4" sheet (with fund-holdings data) right after the
# "总计" summary tab, and roll the new quarter's totals into "总计".
#
# NOTE: reading `.Value` back from a Range in this host returns the
# property descriptor instead of the cell's contents, so reads use
# `.Value2` throughout; writes still use plain `.Value`.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q4" worksheet by duplicating the existing
#    "2022-Q3" sheet (2nd tab) so it inherits the same header/row
#    styling, then drop the copy in right before that sheet.
# ------------------------------------------------------------------
$templateWs = $wb.Worksheets.Item(2)
$templateWs.Copy($templateWs)
$newWs = $wb.Worksheets.Item(2)
$newWs.Name = "2022-Q4"

# The template only had 4 rows (header + 3 funds); we need 7 (header +
# 6 funds). Extend the data-row formatting down to row 7 by copying
# row 2's look onto rows 5-7.
$newWs.Range("A2:H2").Copy()
$newWs.Range("A5:H7").PasteSpecial(-4122)

# Columns B (fund code) and D:G (text-formatted numbers, e.g. "11.19")
# must stay text, not be coerced into numbers (which would also eat
# leading zeroes in the fund codes) - force Text format before writing.
$newWs.Range("B2:B7").NumberFormat = "@"
$newWs.Range("D2:G7").NumberFormat = "@"

# -- Row 2 --
$newWs.Cells.Item(2,1).Value = 0
$newWs.Cells.Item(2,2).Value = "159636"
$newWs.Cells.Item(2,3).Value = "工银瑞信国证港股通科技ETF"
$newWs.Cells.Item(2,4).Value = "11.19"
$newWs.Cells.Item(2,5).Value = "98.28"
$newWs.Cells.Item(2,6).Value = "2.40"
$newWs.Cells.Item(2,7).Value = "0.2686"
$newWs.Cells.Item(2,8).Value = 10

# -- Row 3 --
$newWs.Cells.Item(3,1).Value = 1
$newWs.Cells.Item(3,2).Value = "501021"
$newWs.Cells.Item(3,3).Value = "华宝标普香港上市中国中小盘指数（LOF）A"
$newWs.Cells.Item(3,4).Value = "4.74"
$newWs.Cells.Item(3,5).Value = "94.57"
$newWs.Cells.Item(3,6).Value = "1.73"
$newWs.Cells.Item(3,7).Value = "0.0820"
$newWs.Cells.Item(3,8).Value = 5

# -- Row 4 --
$newWs.Cells.Item(4,1).Value = 2
$newWs.Cells.Item(4,2).Value = "513160"
$newWs.Cells.Item(4,3).Value = "银华恒生港股通中国科技ETF"
$newWs.Cells.Item(4,4).Value = "0.86"
$newWs.Cells.Item(4,5).Value = "94.26"
$newWs.Cells.Item(4,6).Value = "5.38"
$newWs.Cells.Item(4,7).Value = "0.0463"
$newWs.Cells.Item(4,8).Value = 7

# -- Row 5 --
$newWs.Cells.Item(5,1).Value = 3
$newWs.Cells.Item(5,2).Value = "004266"
$newWs.Cells.Item(5,3).Value = "招商沪港深科技创新主题精选灵活配置混合A"
$newWs.Cells.Item(5,4).Value = "0.93"
$newWs.Cells.Item(5,5).Value = "90.63"
$newWs.Cells.Item(5,6).Value = "2.95"
$newWs.Cells.Item(5,7).Value = "0.0274"
$newWs.Cells.Item(5,8).Value = 8

# -- Row 6 --
$newWs.Cells.Item(6,1).Value = 4
$newWs.Cells.Item(6,2).Value = "006127"
$newWs.Cells.Item(6,3).Value = "华宝标普香港上市中国中小盘指数（LOF）C"
$newWs.Cells.Item(6,4).Value = "0.45"
$newWs.Cells.Item(6,5).Value = "94.57"
$newWs.Cells.Item(6,6).Value = "1.73"
$newWs.Cells.Item(6,7).Value = "0.0078"
$newWs.Cells.Item(6,8).Value = 5

# -- Row 7 --
$newWs.Cells.Item(7,1).Value = 5
$newWs.Cells.Item(7,2).Value = "010754"
$newWs.Cells.Item(7,3).Value = "招商沪港深科技创新主题精选灵活配置混合C"
$newWs.Cells.Item(7,4).Value = "0.25"
$newWs.Cells.Item(7,5).Value = "90.63"
$newWs.Cells.Item(7,6).Value = "2.95"
$newWs.Cells.Item(7,7).Value = "0.0074"
$newWs.Cells.Item(7,8).Value = 8

# Restore the plain/General display format (the cells stay text-typed
# because their stored content is non-numeric to the engine once
# written under Text format) by re-pasting the format from column H,
# which was never touched and is still General.
$newWs.Range("H2").Copy()
$newWs.Range("B2:B7").PasteSpecial(-4122)
$newWs.Range("D2:G7").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert the new 2022-Q4 totals
#    at row 2 and push the existing quarters down by one row. The "A"
#    column is just a 0-based running index, so it is left as-is and
#    only extended for the newly appended last row.
# ------------------------------------------------------------------
$summaryWs = $wb.Worksheets.Item("总计")

for ($r = 7; $r -ge 2; $r--) {
  $dst = $r + 1
  $summaryWs.Cells.Item($dst,2).Value = $summaryWs.Cells.Item($r,2).Value2
  $summaryWs.Cells.Item($dst,3).Value = $summaryWs.Cells.Item($r,3).Value2
  $summaryWs.Cells.Item($dst,4).Value = $summaryWs.Cells.Item($r,4).Value2
}

# New row 8's "A" index cell, styled like the rest of the column.
$summaryWs.Range("A7").Copy()
$summaryWs.Range("A8").PasteSpecial(-4122)
$summaryWs.Cells.Item(8,1).Value = 6

# Write the new 2022-Q4 totals into row 2.
$summaryWs.Cells.Item(2,2).Value = "2022-Q4"
$summaryWs.Cells.Item(2,3).Value = 6
$summaryWs.Cells.Item(2,4).Value = 0.44
